$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh inserts a new week's record (date 2022-12-23 / serial 44918)
# at the top of this variety's data block (rows 663-664), pushing every
# existing record down by two rows. The two oldest records that fall off
# the bottom of the original range (old rows 733-734) land as new rows
# 735-736 at the end of the sheet - a plain insert of two rows at 663
# reproduces exactly that shift.
$ws.Rows("663:664").Insert()

# Row 663 - "Primera" quality, new week
$ws.Cells.Item(663, 1).Value = 8
$ws.Cells.Item(663, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(663, 3).Value = "Coquimbo"
$ws.Cells.Item(663, 4).Value = 44918
$ws.Cells.Item(663, 5).Value = 4
$ws.Cells.Item(663, 6).Value = 100112043
$ws.Cells.Item(663, 7).Value = "Pepino ensalada"
$ws.Cells.Item(663, 8).Value = "Sin especificar"
$ws.Cells.Item(663, 9).Value = "Primera"
$ws.Cells.Item(663, 10).Value = 760
$ws.Cells.Item(663, 11).Value = 14000
$ws.Cells.Item(663, 12).Value = 15000
$ws.Cells.Item(663, 13).Value = 14500
$ws.Cells.Item(663, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(663, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(663, 16).Value = 242
$ws.Cells.Item(663, 17).Value = 60
$ws.Cells.Item(663, 18).Value = "Hortaliza"

# Row 664 - "Segunda" quality, new week
$ws.Cells.Item(664, 1).Value = 8
$ws.Cells.Item(664, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(664, 3).Value = "Coquimbo"
$ws.Cells.Item(664, 4).Value = 44918
$ws.Cells.Item(664, 5).Value = 4
$ws.Cells.Item(664, 6).Value = 100112043
$ws.Cells.Item(664, 7).Value = "Pepino ensalada"
$ws.Cells.Item(664, 8).Value = "Sin especificar"
$ws.Cells.Item(664, 9).Value = "Segunda"
$ws.Cells.Item(664, 10).Value = 400
$ws.Cells.Item(664, 11).Value = 9000
$ws.Cells.Item(664, 12).Value = 10000
$ws.Cells.Item(664, 13).Value = 9500
$ws.Cells.Item(664, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(664, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(664, 16).Value = 119
$ws.Cells.Item(664, 17).Value = 80
$ws.Cells.Item(664, 18).Value = "Hortaliza"
